# Earnings Season Q125 part 1 - add "2024 EV/R" column (EV / 2024 Revenue)
# Inserted as a new column before the existing "2024 EV" (X) -> old (Y) column,
# i.e. becomes the new column Y, pushing everything from Y onward one column right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# 1. Insert a new blank column at position Y (25th column). Excel shifts all
#    formulas/references in columns Y.. to the right automatically.
$ws.Columns.Item(25).Insert()

# 2. Header cell (row 2) - picks up the "2024 EV" header style automatically
#    from its left neighbour (X2); just set the label.
$ws.Range("Y2").Value = "2024 EV/R"

# 3. Data rows: EV / Revenue multiple = (2024 EV) column G over (2024 Revenue) column X
$dataRows = @(4,5,6,7,8,10,11,13,14,18,21)
foreach ($r in $dataRows) {
    $cell = $ws.Range("Y$r")
    $cell.NumberFormat = "0.0\x"
    $cell.Formula = "=G$r/X$r"
}

# 4. Summary row (row 3): extend the TRIMMEAN shared-formula block that used
#    to stop at W3 so it also covers the new column.
$ws.Range("Z3").Copy()
$ws.Range("Y3").PasteSpecial(-4122)   # xlPasteFormats - copy Z3's look (old Y3 style) first
$ws.Range("Y3").NumberFormat = "0.0\x"
$ws.Range("Y3").Formula = "=TRIMMEAN(Y4:Y1048576,80%)"

# 5. Restore the selection to B5 (matches the saved view state of the edit).
$ws.Range("B5").Select()

$wb.Save()
